$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the "BQ" / "1" / ":" / " " runs into a single run reading
#    "BQ1: " (same Tahoma/222222/shd formatting as the original runs).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("BQ1: ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "BQ1: ", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) The "_GoBack" bookmark moves from the end of the first paragraph
#    to wrap around the second picture (Picture 4). Remove it from its
#    old location first.
# ---------------------------------------------------------------------
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

# ---------------------------------------------------------------------
# 3) Picture 1 (first inline picture) gets cropped and resized:
#      srcRect l=1795 r=1410 b=7359 (1000ths of a percent)
#      displayed size cx=5753100 cy=4404360 EMU (453pt x 346.8pt)
#    plus the line/extLst/bwMode/rotWithShape bookkeeping Word writes
#    when a crop is applied through the UI.
# ---------------------------------------------------------------------
$pic1 = $d.InlineShapes.Item(1)
$pic1Range = $d.Range($pic1.Range.Start, $pic1.Range.End)

$drawing1 = @'
<wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="2CD917F6" wp14:editId="68FF57EA"><wp:extent cx="5753100" cy="4404360"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1" name=""/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill rotWithShape="1"><a:blip r:embed="rId10"/><a:srcRect l="1795" r="1410" b="7359"/><a:stretch/></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5753100" cy="4404360"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:ln><a:noFill/></a:ln><a:extLst><a:ext uri="{53640926-AAD7-44D8-BBD7-CCE9431645EC}"><a14:shadowObscured xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"/></a:ext></a:extLst></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline>
'@

$pkg1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="272F3D84" w14:textId="4102F7B8" w:rsidR="0070216E" w:rsidRDefault="008779AD"><w:r><w:rPr><w:noProof/></w:rPr><w:drawing>' + $drawing1 + '</w:drawing></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$pic1Range.InsertXML($pkg1) | Out-Null

# ---------------------------------------------------------------------
# 4) Wrap the "_GoBack" bookmark around Picture 4 (the picture in the
#    next paragraph) — bookmarkStart right before its run, bookmarkEnd
#    right after it.
# ---------------------------------------------------------------------
$pic2 = $d.InlineShapes.Item(2)
$pic2Range = $d.Range($pic2.Range.Start, $pic2.Range.End)

$run2 = @'
<w:r><w:rPr><w:noProof/></w:rPr><w:lastRenderedPageBreak/><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="546AD1D4" wp14:editId="2C1142B7"><wp:extent cx="5943600" cy="4754245"/><wp:effectExtent l="0" t="0" r="0" b="8255"/><wp:docPr id="4" name="Picture 4"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1" name=""/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId11"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5943600" cy="4754245"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r>
'@

$content2 = '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + $run2 + '<w:bookmarkEnd w:id="0"/>'

$pkg2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="2390224C" w14:textId="7713FA00" w:rsidR="008779AD" w:rsidRDefault="008779AD">' + $content2 + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$pic2Range.InsertXML($pkg2) | Out-Null

Write-Host "Edit complete"
